# Apply the Player.xlsx "Property" sheet edits described in the commit:
#   "modify queue lock, modify the way of saving player'data"
#
# 1) Flip the "Save" (column E) boolean flag from TRUE to FALSE for the
#    MAXHP/MAXMP-style stat rows 44..67 (the "queue lock" columns).
# 2) Drop the highlighted/"new row" formatting on rows 76 and 77 (GameID /
#    GateID) now that they are no longer pending additions, restoring them
#    to plain default formatting (A76 keeps the text number format).
# 3) Move the saved cursor selection to H78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Column E (Save) TRUE -> FALSE for rows 44 through 67 -----------
$ws.Range("E44:E67").Value = $false

# --- 2) Clear the highlighted formatting from rows 76 and 77 -----------
$plainCells = @("B76","G76","H76","I76","J76","A77","B77","G77","H77","I77","J77")
foreach ($addr in $plainCells) {
    $ws.Range($addr).ClearFormats()
}

# A76 keeps the "@" text number format (style index 1) but loses the
# highlight font/fill, so clear then re-apply just the number format.
$ws.Range("A76").ClearFormats()
$ws.Range("A76").NumberFormat = "@"

# --- 3) Update the remembered selection to H78 --------------------------
$ws.Range("H78").Select()
